$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A51").Value = "Ashirboyeva Shoxsanam Shoikromovna"
$ws.Range("B51").Value = "Maktabgacha talim tashkiloti direktori"
$ws.Range("C51").Value = "AD5815225"
$ws.Range("E51").Value = "Toshkent shahri"
$ws.Range("F51").Value = "Mirzo Ulugʻbek tumani"

# D51, G51, H51 hold values that *look* numeric/date ("217", "+998909794434",
# "03-11-2024") but must be stored as literal text, matching the rest of the
# sheet (every other cell uses inline/shared strings, not numbers/dates).
# A direct .Value assignment gets auto-coerced into a Number/Date. To avoid
# that while still ending up with the default (unstyled) cell format, stage
# the text in a scratch cell that is explicitly formatted as Text, copy it,
# and paste-special (values only) into the destination: this carries over
# the text *value* but not the scratch cell's number format, so the
# destination keeps style 0 exactly like the surrounding cells.
$scratch = $ws.Range("Z1")

$scratch.NumberFormat = "@"
$scratch.Value = "217"
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$scratch.Clear()

$scratch.NumberFormat = "@"
$scratch.Value = "+998909794434"
$scratch.Copy()
$ws.Range("G51").PasteSpecial(-4163)
$scratch.Clear()

$scratch.NumberFormat = "@"
$scratch.Value = "03-11-2024"
$scratch.Copy()
$ws.Range("H51").PasteSpecial(-4163)
$scratch.Clear()
